$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("B1").Value = "Variável"
$ws.Range("C1").Value = "Valor"
$ws.Range("D1").Value = "Colocação"

# Copy style from B1 (already bold/centered) to the new header cells
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$variavel = "Diferença 2023/03 - 2022/03"

# Row 2: Acre
$ws.Range("A2").Value = "Acre"
$ws.Range("B2").Value = $variavel
$ws.Range("C2").Value = 0.03823529411764703
$ws.Range("D2").Value = "1º"

# Row 3: Maranhão
$ws.Range("A3").Value = "Maranhão"
$ws.Range("B3").Value = $variavel
$ws.Range("C3").Value = 0.02872575831950319
$ws.Range("D3").Value = "2º"

# Row 4: Sergipe
$ws.Range("A4").Value = "Sergipe"
$ws.Range("B4").Value = $variavel
$ws.Range("C4").Value = 0.02332956801533514
$ws.Range("D4").Value = "3º"

# Row 5: Distrito Federal
$ws.Range("A5").Value = "Distrito Federal"
$ws.Range("B5").Value = $variavel
$ws.Range("C5").Value = 0.02073694380841284
$ws.Range("D5").Value = "4º"

# Row 6: Bahia
$ws.Range("A6").Value = "Bahia"
$ws.Range("B6").Value = $variavel
$ws.Range("C6").Value = 0.01795666289274156
$ws.Range("D6").Value = "5º"

# Row 7: Espírito Santo
$ws.Range("A7").Value = "Espírito Santo"
$ws.Range("B7").Value = $variavel
$ws.Range("C7").Value = 0.01738161528857163
$ws.Range("D7").Value = "6º"

# Row 8: Nordeste (no Colocação)
$ws.Range("A8").Value = "Nordeste"
$ws.Range("B8").Value = $variavel
$ws.Range("C8").Value = 0.01121801151938073
$ws.Range("D8").Value = ""

# Row 9: Brasil (no Colocação)
$ws.Range("A9").Value = "Brasil"
$ws.Range("B9").Value = $variavel
$ws.Range("C9").Value = 0.01011494672972024
$ws.Range("D9").Value = ""
